$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Process-user matrix")
$ws2 = $wb.Worksheets.Item("Entity - process matrix")
$ws3 = $wb.Worksheets.Item("Entity-user matrix")

# --- Sheet 1: Process-user matrix -------------------------------------
$ws1.Activate()
$ws1.Range("D8").Select()

# --- Sheet 2: Entity - process matrix ----------------------------------
$ws2.Activate()
$ws2.Range("L5").Value = "I,U,D,S"
$ws2.Range("L7").Value = "I,U,D,S"
$ws2.Range("L8").Value = "I,U,D,S"
$ws2.Range("L9").Value = "I,U,D,S"
$ws2.Range("H11").Value = "I,U,D,S"
$ws2.Range("F12").Value = "I,U"
$ws2.Range("H15").Value = "I,U,D,S"
$ws2.Range("I15").Value = "I,U,D,S"
$ws2.Range("I16").Value = "I,U,D,S"
$ws2.Range("H15").Select()

# --- Sheet 3: Entity-user matrix ---------------------------------------
$ws3.Activate()
$ws3.Range("D5").Value = "I,U,D,S"
$ws3.Range("D7").Value = "I,U,D,S"
$ws3.Range("D8").Value = "I,U,D,S"
$ws3.Range("D9").Value = "I,U,D,S"
$ws3.Range("E11").Value = "I,U,D,S"
$ws3.Range("B12").Value = "I,U"
$ws3.Range("C15").Value = "I,U,D,S"
$ws3.Range("C16").Value = "I,U,D,S"
$ws3.Range("E16").Value = "I,U,D,S"
$ws3.Range("C20").Select()
